$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), matching the style of the existing
# header row (e.g. G1) which is bold / bordered / centered.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Fill the new column's data rows with 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
